$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 234, shifting existing rows 234:242 down to 235:243
$ws.Rows.Item(234).Insert()

# Populate the newly inserted row 234 with the new record
$ws.Range("A234").Value = 4
$ws.Range("B234").Value = 'Feria Lagunitas de Puerto Montt'
$ws.Range("C234").Value = 'Los Lagos'
$ws.Range("D234").Value = 44753
$ws.Range("E234").Value = 10
$ws.Range("F234").Value = 'Fruta'
$ws.Range("G234").Value = 100102
$ws.Range("H234").Value = 'Cítricos'
$ws.Range("I234").Value = 100102004
$ws.Range("J234").Value = 'Mandarina'
$ws.Range("K234").Value = 'Clemenuless'
$ws.Range("L234").Value = 'Primera'
$ws.Range("M234").Value = 400
$ws.Range("N234").Value = 9000
$ws.Range("O234").Value = 9000
$ws.Range("P234").Value = 9000
$ws.Range("Q234").Value = '$/bandeja 10 kilos'
$ws.Range("R234").Value = 'Provincia de Limarí'
$ws.Range("S234").Value = 900
$ws.Range("T234").Value = 10

# Match the date format used by the rest of column D
$ws.Range("D234").NumberFormat = $ws.Range("D235").NumberFormat
